$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.045.35"
$ws.Range("E2").Value = "  +1.43%  "

# Row 3
$ws.Range("D3").Value = "2.060.99"
$ws.Range("E3").Value = "  -1.85%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'249.56"
$ws.Range("E5").Value = "  -1.45%  "

# Row 6
$ws.Range("D6").Value = "'0.671"
$ws.Range("E6").Value = "  +2.16%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "'54.89"
$ws.Range("E8").Value = "  +15.39%  "

# Row 9
$ws.Range("D9").Value = "'60.95"
$ws.Range("E9").Value = "  +1.62%  "

# Row 10
$ws.Range("D10").Value = "'0.384"
$ws.Range("E10").Value = "  +1.05%  "

# Row 11
$ws.Range("D11").Value = "'0.0786"
$ws.Range("E11").Value = "  +6.00%  "

# Row 12
$ws.Range("E12").Value = "  +5.87%  "

# Row 13
$ws.Range("D13").Value = "'15.08"
$ws.Range("E13").Value = "  +3.81%  "

# Row 14
$ws.Range("D14").Value = "2.362.84"
$ws.Range("E14").Value = "  -1.77%  "

# Row 15
$ws.Range("D15").Value = "'0.817"
$ws.Range("E15").Value = "  -2.30%  "

# Row 16
$ws.Range("D16").Value = "'5.33"
$ws.Range("E16").Value = "  +4.92%  "

# Row 17
$ws.Range("D17").Value = "2.063.92"
$ws.Range("E17").Value = "  -1.67%  "

# Row 18
$ws.Range("D18").Value = "37.038.51"
$ws.Range("E18").Value = "  +1.44%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0927"
$ws.Range("E19").Value = "  +11.65%  "

# Row 20
$ws.Range("D20").Value = "'73.51"
$ws.Range("E20").Value = "  +1.07%  "

# Row 21
$ws.Range("D21").Value = "'14.25"
$ws.Range("E21").Value = "  +7.72%  "

# Row 22
$ws.Range("E22").Value = "  +2.87%  "

# Row 23
$ws.Range("D23").Value = "'237.57"
$ws.Range("E23").Value = "  -1.12%  "

# Row 24
$ws.Range("E24").Value = "  +0.00%  "

# Row 25
$ws.Range("E25").Value = "  -3.41%  "

# Row 26
$ws.Range("D26").Value = "'170.10"
$ws.Range("E26").Value = "  -0.69%  "

# Row 27
$ws.Range("D27").Value = "'9.03"
$ws.Range("E27").Value = "  -1.30%  "

# Row 28
$ws.Range("D28").Value = "'20.13"
$ws.Range("E28").Value = "  -5.78%  "

# Row 29
$ws.Range("E29").Value = "  +0.25%  "

# Row 30
$ws.Range("E30").Value = "  +1.44%  "

# Row 31
$ws.Range("E31").Value = "  +2.58%  "

# Row 32
$ws.Range("E32").Value = "  +7.75%  "

# Row 33
$ws.Range("D33").Value = "'0.0629"
$ws.Range("E33").Value = "  +2.18%  "

# Row 34
$ws.Range("D34").Value = "'4.41"
$ws.Range("E34").Value = "  +7.90%  "

# Row 35
$ws.Range("D35").Value = "'0.0891"
$ws.Range("E35").Value = "  -2.54%  "

# Row 37
$ws.Range("E37").Value = "  -6.22%  "

# Row 38
$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -4.57%  "

# Row 39
$ws.Range("E39").Value = "  +0.34%  "

# Row 40
$ws.Range("E40").Value = "  +22.35%  "

# Row 41
$ws.Range("E41").Value = "  +0.81%  "

# Row 42
$ws.Range("D42").Value = "'17.73"
$ws.Range("E42").Value = "  +11.19%  "

# Row 43
$ws.Range("E43").Value = "  -1.72%  "

# Row 44
$ws.Range("D44").Value = "'96.95"
$ws.Range("E44").Value = "  -1.14%  "

# Row 45
$ws.Range("E45").Value = "  +1.18%  "

# Row 46
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.15"
$ws.Range("E46").Value = "  +46.57%  "

# Row 47
$ws.Range("B47").Value = "Gas"
$ws.Range("C47").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D47").Value = "'13.76"
$ws.Range("E47").Value = "  -52.33%  "

# Row 48
$ws.Range("E48").Value = "  +6.76%  "

# Row 49
$ws.Range("D49").Value = "1.296.30"
$ws.Range("E49").Value = "  -2.70%  "

# Row 50
$ws.Range("E50").Value = "  +2.29%  "

# Row 51
$ws.Range("D51").Value = "'4.13"
$ws.Range("E51").Value = "  +8.35%  "

